$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing A:E data to B:F
$ws.Columns.Item(1).Insert()

# Populate new column A: header "Metodo" plus the 8 method names
$metodoValues = @("Metodo", "SMARTER", "Fuzzy", "TOPSIS", "GRA", "CODAS", "MABAC", "VIKOR", "PROMETHEE II")
for ($i = 0; $i -lt $metodoValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $metodoValues[$i]
}

# Rename the shifted headers in row 1 (columns B:F) from Var1_x to the real metric names
$headerValues = @("Rx", "Ry", "CL", "Entropia", "SSIM")
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(1, $col).Value = $headerValues[$i]
}

# Resize columns to fit their new contents
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(3).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666
